# Update "想去人数" (want-to-go count) figures in column F for rows 4-12 and 14
# on both the "展览" (sheet1) and "全部类型" (sheet4) worksheets, matching the
# refreshed data snapshot ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" -- row => new F value
$exhibition = $wb.Worksheets.Item("展览")
$exhibition.Range("F4").Value = 70
$exhibition.Range("F5").Value = 535
$exhibition.Range("F6").Value = 7389
$exhibition.Range("F7").Value = 467
$exhibition.Range("F8").Value = 182
$exhibition.Range("F9").Value = 1066
$exhibition.Range("F10").Value = 504
$exhibition.Range("F11").Value = 15
$exhibition.Range("F12").Value = 158
$exhibition.Range("F14").Value = 683

# Sheet "全部类型" -- same rows, note F6 differs slightly (7390 vs 7389)
$allTypes = $wb.Worksheets.Item("全部类型")
$allTypes.Range("F4").Value = 70
$allTypes.Range("F5").Value = 535
$allTypes.Range("F6").Value = 7390
$allTypes.Range("F7").Value = 467
$allTypes.Range("F8").Value = 182
$allTypes.Range("F9").Value = 1066
$allTypes.Range("F10").Value = 504
$allTypes.Range("F11").Value = 15
$allTypes.Range("F12").Value = 158
$allTypes.Range("F14").Value = 683
